$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3075
$ws.Range("I2").Value = 2066.3333
$ws.Range("K2").Value = 2066.3333
$ws.Range("M2").Value = -1953.3333

$ws.Range("H82").Value = 570.1818
$ws.Range("I82").Value = 607.3
$ws.Range("J82").Value = 199
$ws.Range("K82").Value = 1821.9
$ws.Range("L82").Value = 597
$ws.Range("M82").Value = -1415.9
$ws.Range("N82").Value = -1409

$ws.Range("H85").Value = 570.1818
$ws.Range("I85").Value = 607.3
$ws.Range("J85").Value = 199
$ws.Range("K85").Value = 1821.9
$ws.Range("L85").Value = 597
$ws.Range("M85").Value = -417.8999999999999
$ws.Range("N85").Value = -3405

$ws.Range("H116").Value = 1889.5151
$ws.Range("I116").Value = 1721.3182
$ws.Range("J116").Value = 2225.9092
$ws.Range("K116").Value = 1721.3182
$ws.Range("L116").Value = 2225.9092
$ws.Range("M116").Value = 1720.6818
$ws.Range("N116").Value = -9109.9092

$ws.Range("H132").Value = 1176.1578
$ws.Range("I132").Value = 1152.0944
$ws.Range("J132").Value = 1495
$ws.Range("K132").Value = 3456.2832
$ws.Range("L132").Value = 4485
$ws.Range("M132").Value = -926.2831999999999
$ws.Range("N132").Value = -9545

$ws.Range("H135").Value = 38462256
$ws.Range("I135").Value = 13889615
$ws.Range("J135").Value = 333333920
$ws.Range("K135").Value = 125006535
$ws.Range("L135").Value = 3000005280
$ws.Range("M135").Value = -125004000
$ws.Range("N135").Value = -3000010350

$ws.Range("H137").Value = 1840.1818
$ws.Range("I137").Value = 1432.3611
$ws.Range("J137").Value = 3675.375
$ws.Range("K137").Value = 4297.0833
$ws.Range("L137").Value = 11026.125
$ws.Range("M137").Value = -1747.0833
$ws.Range("N137").Value = -16126.125

$ws.Range("H138").Value = 2516.9644
$ws.Range("I138").Value = 1150.0817
$ws.Range("J138").Value = 4430.6
$ws.Range("K138").Value = 3450.2451
$ws.Range("L138").Value = 13291.8
$ws.Range("M138").Value = 1689.7549
$ws.Range("N138").Value = -23571.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18394.324
$ws.Range("I32").Value = 19494.82
$ws.Range("J32").Value = 11681.3
$ws.Range("K32").Value = 19494.82
$ws.Range("L32").Value = 11681.3
$ws.Range("M32").Value = -19207.82
$ws.Range("N32").Value = -12255.3

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws.Range("H132").Value = 2565.772
$ws.Range("I132").Value = 2096.4524
$ws.Range("J132").Value = 3879.8667
$ws.Range("K132").Value = 6289.3572
$ws.Range("L132").Value = 11639.6001
$ws.Range("M132").Value = -3759.3572
$ws.Range("N132").Value = -16699.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 35000
$ws.Range("J51").Value = 35000
$ws.Range("L51").Value = 35000
$ws.Range("N51").Value = -35982

$ws.Range("H55").Value = 60500
$ws.Range("J55").Value = 60500
$ws.Range("L55").Value = 60500
$ws.Range("N55").Value = -61046

$ws.Range("H99").Value = 1082.8334
$ws.Range("I99").Value = 1024.5
$ws.Range("J99").Value = 1199.5
$ws.Range("K99").Value = 1024.5
$ws.Range("L99").Value = 1199.5
$ws.Range("M99").Value = 473.5
$ws.Range("N99").Value = -4195.5

$ws.Range("H105").Value = 923677.3
$ws.Range("I105").Value = 1306451.1
$ws.Range("J105").Value = 5020.1
$ws.Range("K105").Value = 1306451.1
$ws.Range("L105").Value = 5020.1
$ws.Range("M105").Value = -1304704.1
$ws.Range("N105").Value = -8514.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 133.88889
$ws.Range("I22").Value = 141.4
$ws.Range("J22").Value = 124.5
$ws.Range("K22").Value = 141.4
$ws.Range("L22").Value = 124.5
$ws.Range("M22").Value = 208.6
$ws.Range("N22").Value = -824.5

$ws.Range("H31").Value = 2062.3157
$ws.Range("I31").Value = 1466.8667
$ws.Range("J31").Value = 4295.25
$ws.Range("K31").Value = 1466.8667
$ws.Range("L31").Value = 4295.25
$ws.Range("M31").Value = -1171.8667
$ws.Range("N31").Value = -4885.25

$ws.Range("H34").Value = 2062.3157
$ws.Range("I34").Value = 1466.8667
$ws.Range("J34").Value = 4295.25
$ws.Range("K34").Value = 1466.8667
$ws.Range("L34").Value = 4295.25
$ws.Range("M34").Value = -1264.8667
$ws.Range("N34").Value = -4699.25

$ws.Range("H58").Value = 1517080.4
$ws.Range("I58").Value = 2115392
$ws.Range("J58").Value = 3703.5293
$ws.Range("K58").Value = 2115392
$ws.Range("L58").Value = 3703.5293
$ws.Range("M58").Value = -2115189
$ws.Range("N58").Value = -4109.5293

$ws.Range("H94").Value = 1235.5385
$ws.Range("I94").Value = 1005
$ws.Range("J94").Value = 1277.4546
$ws.Range("K94").Value = 1005
$ws.Range("L94").Value = 1277.4546
$ws.Range("M94").Value = -554
$ws.Range("N94").Value = -2179.4546

$ws.Range("H122").Value = 5886
$ws.Range("I122").Value = 5793.6665
$ws.Range("J122").Value = 6440
$ws.Range("K122").Value = 17380.9995
$ws.Range("L122").Value = 19320
$ws.Range("M122").Value = -14930.9995
$ws.Range("N122").Value = -24220

$ws.Range("H132").Value = 2707.2683
$ws.Range("I132").Value = 2944.7666
$ws.Range("J132").Value = 2059.5454
$ws.Range("K132").Value = 8834.299800000001
$ws.Range("L132").Value = 6178.6362
$ws.Range("M132").Value = -6304.299800000001
$ws.Range("N132").Value = -11238.6362

$ws.Range("H134").Value = 2017.1428
$ws.Range("I134").Value = 1266.6666
$ws.Range("J134").Value = 3259.3103
$ws.Range("K134").Value = 3799.9998
$ws.Range("L134").Value = 9777.930899999999
$ws.Range("M134").Value = -1264.9998
$ws.Range("N134").Value = -14847.9309

$ws.Range("H136").Value = 1517080.4
$ws.Range("I136").Value = 2115392
$ws.Range("J136").Value = 3703.5293
$ws.Range("K136").Value = 6346176
$ws.Range("L136").Value = 11110.5879
$ws.Range("M136").Value = -6343626
$ws.Range("N136").Value = -16210.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 585.6667
$ws.Range("I41").Value = 80
$ws.Range("J41").Value = 838.5
$ws.Range("K41").Value = 240
$ws.Range("L41").Value = 2515.5
$ws.Range("M41").Value = 98
$ws.Range("N41").Value = -3191.5

$ws.Range("H75").Value = 3398.9
$ws.Range("I75").Value = 495
$ws.Range("J75").Value = 4124.875
$ws.Range("K75").Value = 1485
$ws.Range("L75").Value = 12374.625
$ws.Range("M75").Value = -487
$ws.Range("N75").Value = -14370.625

$ws.Range("H78").Value = 3398.9
$ws.Range("I78").Value = 495
$ws.Range("J78").Value = 4124.875
$ws.Range("K78").Value = 4455
$ws.Range("L78").Value = 37123.875
$ws.Range("M78").Value = 537
$ws.Range("N78").Value = -47107.875

$ws.Range("H131").Value = 20297.44
$ws.Range("J131").Value = 23749.572
$ws.Range("L131").Value = 71248.716
$ws.Range("N131").Value = -81328.716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5069.36
$ws.Range("I102").Value = 4819.125
$ws.Range("J102").Value = 5514.222
$ws.Range("K102").Value = 4819.125
$ws.Range("L102").Value = 5514.222
$ws.Range("M102").Value = -3197.125
$ws.Range("N102").Value = -8758.222

$ws.Range("H126").Value = 2772.4783
$ws.Range("I126").Value = 1830.5834
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 5491.7502
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -3021.7502
$ws.Range("N126").Value = -16340

$ws.Range("H132").Value = 6744.6924
$ws.Range("I132").Value = 5023.839
$ws.Range("J132").Value = 13413
$ws.Range("K132").Value = 15071.517
$ws.Range("L132").Value = 40239
$ws.Range("M132").Value = -12541.517
$ws.Range("N132").Value = -45299

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5497.636
$ws.Range("I7").Value = 4605.5454
$ws.Range("J7").Value = 7281.8184
$ws.Range("K7").Value = 4605.5454
$ws.Range("L7").Value = 7281.8184
$ws.Range("M7").Value = -4493.5454
$ws.Range("N7").Value = -7505.8184

$ws.Range("H61").Value = 640224.0600000001
$ws.Range("I61").Value = 16238.667
$ws.Range("J61").Value = 10000005
$ws.Range("K61").Value = 16238.667
$ws.Range("L61").Value = 10000005
$ws.Range("M61").Value = -16036.667
$ws.Range("N61").Value = -10000409

$ws.Range("H113").Value = 640224.0600000001
$ws.Range("I113").Value = 16238.667
$ws.Range("J113").Value = 10000005
$ws.Range("K113").Value = 16238.667
$ws.Range("L113").Value = 10000005
$ws.Range("M113").Value = -14068.667
$ws.Range("N113").Value = -10004345

$ws.Range("H126").Value = 5497.636
$ws.Range("I126").Value = 4605.5454
$ws.Range("J126").Value = 7281.8184
$ws.Range("K126").Value = 13816.6362
$ws.Range("L126").Value = 21845.4552
$ws.Range("M126").Value = -11346.6362
$ws.Range("N126").Value = -26785.4552

$ws.Range("H132").Value = 4235.2593
$ws.Range("I132").Value = 4064.5293
$ws.Range("J132").Value = 4525.5
$ws.Range("K132").Value = 12193.5879
$ws.Range("L132").Value = 13576.5
$ws.Range("M132").Value = -9663.5879
$ws.Range("N132").Value = -18636.5

$ws.Range("H136").Value = 3132.75
$ws.Range("I136").Value = 1646.8909
$ws.Range("J136").Value = 7024.2856
$ws.Range("K136").Value = 4940.6727
$ws.Range("L136").Value = 21072.8568
$ws.Range("M136").Value = -2390.6727
$ws.Range("N136").Value = -26172.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 36579
$ws.Range("J76").Value = 36579
$ws.Range("L76").Value = 36579
$ws.Range("N76").Value = -37209

$ws.Range("H79").Value = 36579
$ws.Range("J79").Value = 36579
$ws.Range("L79").Value = 36579

$ws.Range("H126").Value = 1470.6086
$ws.Range("I126").Value = 1496.2667
$ws.Range("J126").Value = 1422.5
$ws.Range("K126").Value = 4488.800099999999
$ws.Range("L126").Value = 4267.5
$ws.Range("M126").Value = -2018.800099999999
$ws.Range("N126").Value = -9207.5

$ws.Range("H132").Value = 1232.1923
$ws.Range("I132").Value = 582.9459000000001
$ws.Range("K132").Value = 1748.8377
$ws.Range("M132").Value = 781.1623

$ws.Range("H136").Value = 4872.418
$ws.Range("I136").Value = 3142.4043
$ws.Range("J136").Value = 8937.950000000001
$ws.Range("K136").Value = 9427.2129
$ws.Range("L136").Value = 26813.85
$ws.Range("M136").Value = -6877.2129
$ws.Range("N136").Value = -31913.85
